$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.991.45'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.009.32'
$ws.Range('E3').Value = '  +3.10%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.31'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.77'
$ws.Range('E6').Value = '  -3.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').Value = '  -2.13%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.00'
$ws.Range('E10').Value = '  -2.86%  '
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0855'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.96'
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('D14').Value = '3.482.26'
$ws.Range('E14').Value = '  +3.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.59'
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').Value = '2.995.20'
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.01'
$ws.Range('E17').Value = '  +4.04%  '
$ws.Range('D18').Value = '52.020.79'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.38'
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.45'
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.57'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').Value = '0.0₃0970'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.03'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.45'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.72'
$ws.Range('E25').Value = '  -3.19%  '
$ws.Range('E26').Value = '  -2.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.91'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.45'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.108'
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.46'
$ws.Range('E31').Value = '  +8.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.18'
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  +15.17%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '35.82'
$ws.Range('E34').Value = '  -6.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.23'
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0435'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +3.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.79'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.44'
$ws.Range('E41').Value = '  -4.60%  '
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.28'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '124.88'
$ws.Range('E44').Value = '  +3.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.19'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '2.128.81'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -6.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.243'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0337'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.903'
$ws.Range('E51').Value = '  +0.44%  '
